$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3
$ws.Range("F6").Value = -4
$ws.Range("F8").Value = -7
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = -1
$ws.Range("F14").Value = 7
